$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# Sheet 1 (Overview): columns E (zh-cn) and F (de-de), rows 2-4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# Sheet 2 (zh-cn): Status column C, rows 2-4
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2:C4").Value = "In Translation"

# Sheet 3 (de-de): Status column C, rows 2-4
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2:C4").Value = "In Translation"

# --- Narrow the Status-related columns to reflect the shorter text ---
# Overview sheet: columns E and F
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C
$wsZh.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C
$wsDe.Columns.Item(3).ColumnWidth = 12.5
